# Splits the run covering [start,end) away from its neighbors by toggling
# (and then reverting) a character-formatting property on it. This forces
# the document engine to keep it as its own <w:r> instead of silently
# re-coalescing it with adjacent, identically-formatted runs.
function Split-Run {
    param($doc, $start, $end)
    $r1 = $doc.Range($start, $end)
    $r1.Font.Bold = $true
    $r2 = $doc.Range($start, $end)
    $r2.Font.Bold = $false
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Occurrence 1 (first "TID 1: {tid1}" paragraph):
#   "TID 1:" / " {tid1}"  ->  "TI" / "D" / " {key0}" / ":" / " {value0}"
# ---------------------------------------------------------------------
$rng = $d.Content
$found1 = $rng.Find.Execute("TID 1: {tid1}")
$afterOccurrence1 = 0
if ($found1) {
    $s1 = $rng.Start
    $rng.Text = "TID {key0}: {value0}"

    $p1s = $s1;      $p1e = $s1 + 2    # "TI"
    $p2s = $s1 + 2;  $p2e = $s1 + 3    # "D"
    $p3s = $s1 + 3;  $p3e = $s1 + 10   # " {key0}"
    $p4s = $s1 + 10; $p4e = $s1 + 11   # ":"
    # last piece " {value0}" needs no trailing split, it ends the paragraph's run sequence

    Split-Run $d $p1s $p1e
    Split-Run $d $p2s $p2e
    Split-Run $d $p3s $p3e
    Split-Run $d $p4s $p4e

    $afterOccurrence1 = $s1 + 20
}

# ---------------------------------------------------------------------
# Occurrence 2 (second "TID 1: {tid1}" paragraph):
#   "TID 1:" / " {tid1}"  ->  "TID" / " {key0}" / ":" / " {value0}"
# ---------------------------------------------------------------------
$docEnd = $d.Content.End
$rng2 = $d.Range($afterOccurrence1, $docEnd)
$found2 = $rng2.Find.Execute("TID 1: {tid1}")
if ($found2) {
    $s2 = $rng2.Start
    $rng2.Text = "TID {key0}: {value0}"

    $q1s = $s2;      $q1e = $s2 + 3    # "TID"
    $q2s = $s2 + 3;  $q2e = $s2 + 10   # " {key0}"
    $q3s = $s2 + 10; $q3e = $s2 + 11   # ":"
    # last piece " {value0}" needs no trailing split

    Split-Run $d $q1s $q1e
    Split-Run $d $q2s $q2e
    Split-Run $d $q3s $q3e
}

# ---------------------------------------------------------------------
# Occurrence 3 (table cell): simple text swap "{tid1}" -> "{value0}"
# ---------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("{tid1}", $true, $false, $false, $false, $false, $true, 1, $false, "{value0}", 2)
